$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift BGM values: id 1000 becomes "None", existing BGM names shift down one row
$ws.Range("B2").Value = "None"
$ws.Range("B3").Value = "BGM_Boss_01"
$ws.Range("B4").Value = "BGM_Boss_02"
$ws.Range("B5").Value = "BGM_Boss_03"

# Update selection to B2
$ws.Range("B2").Select()
